# reserved_words.xlsx - "bissl zeug gmacht lul"
# Update the reserved-words / symbols table on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column layout: a new narrow "symbol" column is inserted before the
# existing wide "symbol meaning" columns (E keeps its own width now,
# F:G stay at their old shared width). We don't actually insert a real
# column - the data itself already lines up one slot to the right in
# columns E/F/G - we only need to restyle the column width metadata.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.0        # -> ~22.76 (B, new narrow col)
$ws.Columns.Item(5).ColumnWidth = 10.8333333  # -> ~11.71 (E, now its own width)

# ---------------------------------------------------------------------
# Row height tweaks (rows that lost their "tall" content shrink back to
# the normal 13.8 row height used throughout the sheet).
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 13.8

# ---------------------------------------------------------------------
# Cell-level content changes.
# ---------------------------------------------------------------------

# Row 3: the char-literal example moves out of this row entirely.
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4: drop the standalone "char" keyword (column B).
$ws.Range("B4").ClearContents()
$ws.Range("F4").Value = "#234 is an int in hex format"

# Row 5: drop the standalone "short" keyword (column B).
$ws.Range("B5").ClearContents()

# Row 7: drop the standalone "long" keyword (column B).
$ws.Range("B7").ClearContents()

# Row 8: "free" becomes "del".
$ws.Range("A8").Value = "del"

# Row 9: drop the standalone "double" keyword (column B).
$ws.Range("B9").ClearContents()

# Row 11: new char-literal example now lives here.
$ws.Range("E11").Value = "'a" + [char]0x2019
# Setting a value that starts with an apostrophe makes Excel flag the
# cell with a "quote prefix" style; paste the formatting back from a
# neighbouring plain cell in the same column so the style index (s=2)
# stays exactly what it was.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F11").Value = "a is a non-unicode string"

# Row 21: "||" / "amount" become "|a|" / "amount of a".
$ws.Range("E21").Value = "|a|"
$ws.Range("F21").Value = "amount of a"

# Row 22: drop the trailing "foreach" keyword, but keep the (now empty)
# row alive / part of the sheet's used range, matching row height.
$ws.Range("A22").ClearContents()
$ws.Rows.Item(22).RowHeight = 13.8
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: zoom out a bit and move the selection/scroll position.
# ---------------------------------------------------------------------
$ws.Range("E1").Select()
$excel.ActiveWindow.Zoom = 160
try {
    $excel.ActiveWindow.ScrollColumn = 4
} catch {}
